# Move "default_privacy" column from the creators sheet to the users sheet,
# and populate it with initial values (default privacy is now tracked per
# user rather than per creator profile).

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("users")
$wsCreators = $wb.Worksheets.Item("creators")
$wsPosts = $wb.Worksheets.Item("posts")

# --- users sheet: add new column H = default_privacy, with initial values ---
$wsUsers.Activate()

$wsUsers.Range("H1").Value = "default_privacy"
$wsUsers.Range("H1").Font.Bold = $true

$wsUsers.Range("H2").Value = "public"
$wsUsers.Range("H3").Value = "private"
$wsUsers.Range("H4").Value = "private"

$wsUsers.Columns.Item(8).ColumnWidth = 13.166666666666666

$wsUsers.Range("H1").Select() | Out-Null

# --- creators sheet: remove default_privacy column contents ---
$wsCreators.Activate()

$wsCreators.Range("E1:E4").ClearContents()

$wsCreators.Range("E3").Select() | Out-Null

# Restore original active sheet (posts)
$wsPosts.Activate()
